# Update the forecast values (column C) for rows 2 through 23 on the
# active worksheet to reflect the refreshed model output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -43.7768
    3  = -46.8528
    4  = -48.9173
    5  = -56.1118
    6  = -54.793
    7  = -52.7371
    8  = -34.5292
    9  = -46.5199
    10 = -60.7061
    11 = -38.2726
    12 = -54.9382
    13 = -53.0367
    14 = -55.4607
    15 = -61.9898
    16 = -61.3114
    17 = -61.8827
    18 = -66.0309
    19 = -71.0578
    20 = -68.9492
    21 = -67.4239
    22 = -66.6504
    23 = -76.8858
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}
